$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D
$ws.Range("D2").Value = "37.130.96"
$ws.Range("E2").Value = "  +0.37%  "

# Row 3: D
$ws.Range("D3").Value = "2.056.57"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  +0.25%  "

# Row 5: D (text-forced)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("E6").Value = "  -0.47%  "

# Row 7: D (text-forced)
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.92"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.55%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +2.51%  "

# Row 10: D (text-forced)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0796"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.07%  "

$ws.Range("E11").Value = "  +2.11%  "

# Row 12: D (text-forced)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.21%  "

# Row 13: D
$ws.Range("D13").Value = "2.356.25"
$ws.Range("E13").Value = "  +0.19%  "

# Row 14: D (text-forced)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.826"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.62%  "

# Row 15: D (text-forced)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.48%  "

# Row 16: D
$ws.Range("D16").Value = "2.057.40"
$ws.Range("E16").Value = "  +0.18%  "

# Row 17: D (text-forced)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +31.61%  "

# Row 18: D
$ws.Range("D18").Value = "37.129.72"
$ws.Range("E18").Value = "  +0.47%  "

# Row 19: D (text-forced)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.73%  "

# Row 20: D
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  -1.59%  "

# Row 21: D (text-forced)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.03%  "

# Row 22: D (text-forced)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "239.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.39%  "

$ws.Range("E23").Value = "  -0.04%  "

# Row 24: D (text-forced)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25: D (text-forced)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.93%  "

# Row 26: D (text-forced)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.55%  "

# Row 27: D (text-forced)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "

# Row 28: D (text-forced)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.69%  "

$ws.Range("E29").Value = "  +1.21%  "

$ws.Range("E32").Value = "  +0.44%  "

# Row 33: D (text-forced)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.71%  "

# Row 34: D (text-forced)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0895"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.83%  "

$ws.Range("E35").Value = "  +0.03%  "

# Row 36: D (text-forced)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.01%  "

# Row 37: D (text-forced)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.70%  "

$ws.Range("E38").Value = "  +5.98%  "

# Row 39: D (text-forced)
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "

# Row 40: D (text-forced)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +27.25%  "

# Row 41: D (text-forced)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.11%  "

$ws.Range("E43").Value = "  +0.88%  "

$ws.Range("E44").Value = "  +0.72%  "

# Row 45: D (text-forced)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.91%  "

# Row 48: D
$ws.Range("D48").Value = "1.293.64"
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("E49").Value = "  -1.20%  "

$ws.Range("E50").Value = "  +0.70%  "

# Row 51: D
$ws.Range("D51").Value = "2.236.55"
$ws.Range("E51").Value = "  -0.54%  "

# Row 30: swap Filecoin -> ImmutableX
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.46%  "

# Row 31: swap ImmutableX -> Filecoin
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.21%  "

# Row 46: swap FTXToken -> RenderToken
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.28%  "

# Row 47: swap RenderToken -> FTXToken
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.36%  "

